# Update "想去人数" (F column) counts and one Cover image URL (I column)
# on both the "展览" and "全部类型" worksheets, reflecting the latest
# scrape output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value  = 13
$ws1.Range("F3").Value  = 7743
$ws1.Range("F4").Value  = 2711
$ws1.Range("F8").Value  = 590
$ws1.Range("F10").Value = 68
$ws1.Range("F12").Value = 860
$ws1.Range("F13").Value = 3103
$ws1.Range("F14").Value = 196
$ws1.Range("F15").Value = 86
$ws1.Range("F16").Value = 727
$ws1.Range("F17").Value = 752
$ws1.Range("I18").Value = "//i2.hdslb.com/bfs/openplatform/202401/Y7hnq4gB1706517272632.jpeg"
$ws1.Range("F19").Value = 453
$ws1.Range("F21").Value = 226
$ws1.Range("F22").Value = 218
$ws1.Range("F23").Value = 289
$ws1.Range("F25").Value = 126
$ws1.Range("F26").Value = 98
$ws1.Range("F27").Value = 262
$ws1.Range("F28").Value = 9
$ws1.Range("F32").Value = 475
$ws1.Range("F33").Value = 30
$ws1.Range("F34").Value = 19
$ws1.Range("F36").Value = 88

# ---- Sheet "全部类型" ----
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F5").Value  = 13
$ws4.Range("F6").Value  = 7743
$ws4.Range("F7").Value  = 2711
$ws4.Range("F11").Value = 590
$ws4.Range("F13").Value = 68
$ws4.Range("F15").Value = 860
$ws4.Range("F17").Value = 3103
$ws4.Range("F18").Value = 196
$ws4.Range("F19").Value = 86
$ws4.Range("F21").Value = 727
$ws4.Range("F22").Value = 752
$ws4.Range("I24").Value = "//i2.hdslb.com/bfs/openplatform/202401/Y7hnq4gB1706517272632.jpeg"
$ws4.Range("F25").Value = 453
$ws4.Range("F27").Value = 226
$ws4.Range("F28").Value = 218
$ws4.Range("F29").Value = 289
$ws4.Range("F31").Value = 126
$ws4.Range("F32").Value = 98
$ws4.Range("F33").Value = 262
$ws4.Range("F34").Value = 9
$ws4.Range("F38").Value = 475
$ws4.Range("F39").Value = 30
$ws4.Range("F40").Value = 19
$ws4.Range("F42").Value = 88
